$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: grow the AutoFilter range to A1:F79 -------------------------------
# Do this *before* inserting the new row: Range.AutoFilter() (no field/criteria)
# re-filters the sheet's current "used range", so we need the used range to
# still be A1:F79 (i.e. before row 80 exists) for the filter to land on F79.
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("A1:F79").AutoFilter() | Out-Null

# --- Step 2: add the new test-case row (row 80) --------------------------------
# Duplicate row 79 (values + all formatting) into row 80 so the new row inherits
# the same styling (bold-ish ID column, fixed "3,"/"data"/"on"/"off" columns).
$ws.Range("A79:F79").Copy()
$ws.Range("A80:F80").PasteSpecial(-4104) | Out-Null
$ws.Range("A79:F79").Copy()
$ws.Range("A80:F80").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Now overwrite the TestName / ID cells for the new row with the new test data.
$ws.Range("A80").Value = "Domestic_Payments_In_Future_[WEB]"
$ws.Range("B80").Value = "C70835"

# --- Step 3: leave the selection where Excel would after inserting a row -------
$ws.Range("A82").Select() | Out-Null

# --- Step 4: keep the hidden _FilterDatabase defined name in sync with the
# AutoFilter range (Excel normally maintains this automatically).
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$1:`$F`$79"

Write-Host "Added row 80 (Domestic_Payments_In_Future_[WEB] / C70835) and extended AutoFilter to A1:F79"
